$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.181.43"
$ws.Range("D3").Value = "3.424.34"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "414.03"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.22"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  -2.54%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.723"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.74"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.23"
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "3.965.09"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000213"
$ws.Range("E14").Value = "  +3.19%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.140"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.47"
$ws.Range("E16").Value = "  -3.45%  "
$ws.Range("D17").Value = "3.442.40"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.61"
$ws.Range("E18").Value = "  +2.81%  "
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "62.246.53"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "465.48"
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.66"
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.26"
$ws.Range("E23").Value = "  +3.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.48"
$ws.Range("E24").Value = "  +4.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.42"
$ws.Range("E25").Value = "  +17.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.31"
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.00"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.79"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.90"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("E31").Value = "  -3.75%  "
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.112"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.71"
$ws.Range("E34").Value = "  -4.83%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.29"
$ws.Range("E36").Value = "  +9.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0487"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.06"
$ws.Range("E39").Value = "  +4.93%  "
$ws.Range("E40").Value = "  +3.31%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.33"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "145.25"
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("E44").Value = "  +9.57%  "
$ws.Range("E45").Value = "  +5.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.33"
$ws.Range("E46").Value = "  +2.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.54"
$ws.Range("E47").Value = "  +19.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.45"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.37"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").Value = "0.0₃0512"
$ws.Range("E50").Value = "  +23.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.18"
$ws.Range("E51").Value = "  +3.80%  "
